$wb = $excel.ActiveWorkbook

# --- About sheet: insert explanatory note about geothermal -> pumped hydro repurposing ---
$ws1 = $wb.Worksheets.Item("About")
$ws1.Rows("11:12").Insert()
$ws1.Range("A11").Value = "In the India EPS, the geothermal plant type is repurposed as pumped hydro capacity."

# Build the new style (font color black, vertical-centered) on the blank row first,
# then copy/paste the resulting format onto the text row so both cells resolve to the
# same single cell-style entry.
$ws1.Range("A12").Font.Color = 0
$ws1.Range("A12").VerticalAlignment = -4108
$ws1.Range("A12").Copy()
$ws1.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- BPaFF-BITPTaP: flag geothermal (row 10) as providing flexibility ---
$ws2 = $wb.Worksheets.Item("BPaFF-BITPTaP")
$ws2.Range("B10").Value = 1
$ws2.Range("B11").Select()

# --- BPaFF-BDTPTPF: flag geothermal (row 10) as providing flexibility ---
$ws3 = $wb.Worksheets.Item("BPaFF-BDTPTPF")
$ws3.Range("B10").Value = 1
$ws3.Range("B11").Select()

# --- Make the About sheet the active tab/selection ---
$ws1.Activate()
$ws1.Range("B10").Select()
